# Applies the edit described by the diff:
#  - A1 gets the long questions/JSON-style text (as a string)
#  - A1 loses its old bold + bordered + centered style (back to default)
#  - A2 (which used to hold this text) is cleared entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You have been asked to customize the Dynamics 365 Sales application to include a new custom field on the Contact form that will capture the contact's job title. What is the best way to accomplish this task?",
        "ques_type": 2,
        "options": [
            "Create a new extension in the Dynamics 365 Development Environment and add the new field to the Contact form using C# or JavaScript code.",
            "Use the \"Customize this Page\" feature in Dynamics 365 to add the new field to the Contact form.",
            "Utilize the \"Extension Management\" feature in Dynamics 365 to add the new field to the Contact form.",
            "Use the \"Customize Entities\" feature in Dynamics 365 to add the new field to the Contact form."
        ],
        "score": "Create a new extension in the Dynamics 365 Development Environment and add the new field to the Contact form using C# or JavaScript code."
    },
    {
        "title": "You have been asked to create a new custom entity in Dynamics 365 called \"Equipment\" that will be used to track information about your company's equipment. The entity should be related to the existing \"Account\" entity.What is the best way to accomplish this task?",
        "ques_type": 2,
        "options": [
            "Create a new entity in Dynamics 365 using the \"Customize Entities\" feature and manually add the relationships to the \"Account\" entity.",
            "Create a new entity in Dynamics 365 using the \"Customize Entities\" feature, then create a new workflow to automatically create the relationship to the \"Account\" entity.",
            "Create a new entity in Dynamics 365 using the \"Customize Entities\" feature and use the \"Add Existing\" option to add the relationship to the \"Account\" entity.",
            "Create a new extension in the Dynamics 365 Development Environment and add the new entity using C# or JavaScript code, then add the relationship to the \"Account\" entity using code."
        ],
        "score": "Create a new entity in Dynamics 365 using the \"Customize Entities\" feature and use the \"Add Existing\" option to add the relationship to the \"Account\" entity."
    },
    {
        "title": "You have been asked to create a new custom web service in Dynamics 365 that will be used to retrieve information about a specific account. The web service should be secure and should only be accessible to authenticated users. What is the best way to accomplish this task?",
        "ques_type": 2,
        "options": [
            "Use the Dynamics 365 Web Services feature to create a new custom web service and use Basic Authentication to secure the service.",
            "Use the Dynamics 365 Web Services feature to create a new custom web service and use OAuth 2.0 to secure the service.",
            "Use the Dynamics 365 Web Services feature to create a new custom web service and use JSON Web Tokens to secure the service.",
            "Create a new custom web service using C# or JavaScript code and use the Dynamics 365 Web Services feature to secure the service."
        ],
        "score": "Create a new custom web service using C# or JavaScript code and use the Dynamics 365 Web Services feature to secure the service."
    },
    {
        "title": "You have been asked to create a new custom page in Dynamics 365 Business Central that will be used to display a list of open sales orders for a specific customer. The page should include the ability to filter the orders by date range and sort them by order number. What is the best way to accomplish this task?",
        "ques_type": 2,
        "options": [
            "Create a new extension in the Dynamics 365 Business Central Development Environment and add a new page using C/AL code, and include the necessary filters and sorting options using C/AL code.",
            "Use the \"Customize Reports\" feature in Dynamics 365 Business Central to create a new report that includes the necessary filters and sorting options and then create a page to display it.",
            "Utilize the \"Extension Management\" feature in Dynamics 365 Business Central to add a custom page and include the necessary filters and sorting options.",
            "Use the \"Customize this Page\" feature in Dynamics 365 Business Central to create a new custom page and add the necessary filters and sorting options."
        ],
        "score": "Create a new extension in the Dynamics 365 Business Central Development Environment and add a new page using C/AL code, and include the necessary filters and sorting options using C/AL code."
    }
]
'@

# A2 previously held the shared-string text; remove it completely.
$ws.Range("A2").ClearContents()
$ws.Range("A2").ClearFormats()

# A1 previously held 0 with a bold/bordered/centered style; reset the
# formatting back to default before writing the new text into it.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Setting the long text auto-expands the row height; auto-fitting it
# afterwards collapses back to the sheet's default (no explicit <ht>),
# matching the target state where row 1 carries no custom height.
$ws.Rows(1).EntireRow.AutoFit()
